$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 918, shifting the existing data (rows 918:977)
# down to (919:978).
$ws.Rows.Item(918).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A918").Value = 8
$ws.Range("B918").Value = "Terminal La Palmera de La Serena"
$ws.Range("C918").Value = "Coquimbo"
$ws.Range("D918").Value = 44931
$ws.Range("E918").Value = 4
$ws.Range("F918").Value = 100112004
$ws.Range("G918").Value = "Cebolla"
$ws.Range("H918").Value = "Sin especificar"
$ws.Range("I918").Value = "Primera"
$ws.Range("J918").Value = 2000
$ws.Range("K918").Value = 9000
$ws.Range("L918").Value = 10000
$ws.Range("M918").Value = 9500
$ws.Range("N918").Value = "$/malla 18 kilos"
$ws.Range("O918").Value = "Perú"
$ws.Range("P918").Value = 528
$ws.Range("Q918").Value = 18
$ws.Range("R918").Value = "Hortaliza"

# Match the date format used by the other rows in column D.
$ws.Range("D918").NumberFormat = $ws.Range("D919").NumberFormat
